$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so that values
# such as "65.518.74" or "0.110" are not reinterpreted as numbers/dates.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '65.518.74'
$ws.Range('E2').Value = '  -3.40%  '

# Row 3
$ws.Range('D3').Value = '3.489.20'
$ws.Range('E3').Value = '  -1.03%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').Value = '553.15'
$ws.Range('E5').Value = '  -0.79%  '

# Row 6
$ws.Range('D6').Value = '178.86'
$ws.Range('E6').Value = '  -6.83%  '

# Row 7
$ws.Range('D7').Value = '0.641'
$ws.Range('E7').Value = '  +4.59%  '

# Row 8
$ws.Range('E8').Value = '  +0.04%  '

# Row 9
$ws.Range('E9').Value = '  -1.51%  '

# Row 10
$ws.Range('D10').Value = '0.154'
$ws.Range('E10').Value = '  +1.11%  '

# Row 11
$ws.Range('D11').Value = '53.65'
$ws.Range('E11').Value = '  -6.33%  '

# Row 12
$ws.Range('E12').Value = '  -2.30%  '

# Row 13
$ws.Range('D13').Value = '9.21'
$ws.Range('E13').Value = '  -3.28%  '

# Row 14
$ws.Range('D14').Value = '4.042.65'
$ws.Range('E14').Value = '  -0.97%  '

# Row 15
$ws.Range('D15').Value = '3.486.79'
$ws.Range('E15').Value = '  -0.93%  '

# Row 16
$ws.Range('E16').Value = '  +0.04%  '

# Row 17
$ws.Range('D17').Value = '18.41'
$ws.Range('E17').Value = '  +0.09%  '

# Row 18
$ws.Range('D18').Value = '12.16'
$ws.Range('E18').Value = '  +2.13%  '

# Row 19
$ws.Range('D19').Value = '65.472.26'
$ws.Range('E19').Value = '  -4.59%  '

# Row 20
$ws.Range('D20').Value = '0.994'
$ws.Range('E20').Value = '  -1.81%  '

# Row 21
$ws.Range('D21').Value = '412.54'
$ws.Range('E21').Value = '  +0.48%  '

# Row 22
$ws.Range('D22').Value = '4.04'
$ws.Range('E22').Value = '  +1.70%  '

# Row 23
$ws.Range('D23').Value = '85.93'
$ws.Range('E23').Value = '  +1.09%  '

# Row 24
$ws.Range('E24').Value = '  -3.55%  '

# Row 25
$ws.Range('D25').Value = '12.70'
$ws.Range('E25').Value = '  +6.39%  '

# Row 26
$ws.Range('E26').Value = '  -8.26%  '

# Row 27
$ws.Range('E27').Value = '  -2.29%  '

# Row 28
$ws.Range('E28').Value = '  -4.15%  '

# Row 29
$ws.Range('D29').Value = '9.05'
$ws.Range('E29').Value = '  +4.49%  '

# Row 30
$ws.Range('D30').Value = '30.31'
$ws.Range('E30').Value = '  -1.30%  '

# Row 31
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '610.97'
$ws.Range('E31').Value = '  -11.19%  '

# Row 32
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '6.47'
$ws.Range('E32').Value = '  -5.91%  '

# Row 33
$ws.Range('D33').Value = '11.64'
$ws.Range('E33').Value = '  -1.00%  '

# Row 34
$ws.Range('D34').Value = '0.110'
$ws.Range('E34').Value = '  -1.56%  '

# Row 35
$ws.Range('D35').Value = '59.48'
$ws.Range('E35').Value = '  -1.66%  '

# Row 36
$ws.Range('E36').Value = '  +10.37%  '

# Row 37
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.13%  '

# Row 38
$ws.Range('E38').Value = '  -5.65%  '

# Row 39
$ws.Range('D39').Value = '37.09'
$ws.Range('E39').Value = '  -5.65%  '

# Row 40
$ws.Range('D40').Value = '3.361.94'
$ws.Range('E40').Value = '  +9.20%  '

# Row 41
$ws.Range('E41').Value = '  -6.54%  '

# Row 42
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '3.28'
$ws.Range('E42').Value = '  -4.32%  '

# Row 43
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.04%  '

# Row 44
$ws.Range('E44').Value = '  -6.58%  '

# Row 45
$ws.Range('E45').Value = '  -9.21%  '

# Row 46
$ws.Range('E46').Value = '  -1.99%  '

# Row 47
$ws.Range('D47').Value = '3.22'
$ws.Range('E47').Value = '  -0.06%  '

# Row 48
$ws.Range('D48').Value = '2.73'
$ws.Range('E48').Value = '  -0.94%  '

# Row 49
$ws.Range('D49').Value = '0.133'
$ws.Range('E49').Value = '  +1.66%  '

# Row 50
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '137.60'
$ws.Range('E50').Value = '  -1.07%  '

# Row 51
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '8.42'
$ws.Range('E51').Value = '  -10.31%  '

# Restore the original (default) cell style now that the text values are set,
# so the cells keep the workbook's default formatting like before the edit.
$priceVolRange.Style = "Normal"
